$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update numeric values in column B (rows 2-10)
$ws.Range("B2").Value = 141
$ws.Range("B3").Value = 138
$ws.Range("B4").Value = 103
$ws.Range("B5").Value = 101
$ws.Range("B6").Value = 98
$ws.Range("B7").Value = 85
$ws.Range("B8").Value = 78
$ws.Range("B9").Value = 74
$ws.Range("B10").Value = 71

# Swap names in rows 5/6 (CORAS QUISPE JORGE AMERICO and MANOSALVA RUIZ SANDRA KAROLINE)
$ws.Range("A5").Value = "CORAS QUISPE JORGE AMERICO"
$ws.Range("A6").Value = "MANOSALVA RUIZ SANDRA KAROLINE"

# Swap names in rows 9/10 (SAUCEDO CABRERA CARLOS ALEXANDER and CASTREJON TELLO GRECIA)
$ws.Range("A9").Value = "SAUCEDO CABRERA CARLOS ALEXANDER"
$ws.Range("A10").Value = "CASTREJON TELLO GRECIA"

$wb.Save()
